$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.892.45"
$ws.Range("E2").Value = "  +1.33%  "

$ws.Range("D3").Value = "3.461.08"
$ws.Range("E3").Value = "  +3.10%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.02%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "580.79"
$ws.Range("D5").Style = "Normal"

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "148.29"
$ws.Range("D6").Style = "Normal"

$ws.Range("D7").Value = "3.460.12"
$ws.Range("E7").Value = "  +3.12%  "

$ws.Range("E8").Value = "  +0.06%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.474"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.34%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.72"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +3.40%  "

$ws.Range("E11").Value = "  +1.76%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.390"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.41%  "

$ws.Range("D13").Value = "4.052.94"
$ws.Range("E13").Value = "  +3.20%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "28.12"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +8.54%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000176"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.23%  "

$ws.Range("D17").Value = "3.471.79"
$ws.Range("E17").Value = "  +3.55%  "

$ws.Range("D18").Value = "61.936.62"
$ws.Range("E18").Value = "  +1.20%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.35"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +9.28%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.41"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +3.06%  "

$ws.Range("E21").Value = "  +2.44%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "385.31"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.19%  "

$ws.Range("E23").Value = "  +3.29%  "

$ws.Range("D24").Value = "3.594.40"
$ws.Range("E24").Value = "  +3.11%  "

$ws.Range("B25").Value = "Dai"
$ws.Range("C25").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.00"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.24%  "

$ws.Range("B26").Value = "LEO"
$ws.Range("C26").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "5.78"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.06%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "72.64"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.44%  "

$ws.Range("E28").Value = "  -0.96%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.181"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +9.50%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.85"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +5.79%  "

$ws.Range("E31").Value = "  -12.61%  "

$ws.Range("E32").Value = "  -0.40%  "

$ws.Range("E33").Value = "  +1.75%  "

$ws.Range("E34").Value = "  +2.42%  "

$ws.Range("E35").Value = "  +0.00%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "23.98"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.92%  "

$ws.Range("E37").Value = "  +4.47%  "

$ws.Range("E38").Value = "  +0.50%  "

$ws.Range("E39").Value = "  +2.43%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "167.08"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.45%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0790"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +4.47%  "

$ws.Range("B42").Value = "Mantle"
$ws.Range("C42").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.797"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +4.01%  "

$ws.Range("B43").Value = "EnergySwap"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "25.96"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +9.44%  "

$ws.Range("B44").Value = "Stacks"
$ws.Range("C44").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.74"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.45%  "

$ws.Range("B45").Value = "FirstDigitalUSD"
$ws.Range("C45").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.999"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.06%  "

$ws.Range("B46").Value = "OKB"
$ws.Range("C46").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "42.35"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.22%  "

$ws.Range("B47").Value = "Filecoin"
$ws.Range("C47").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.49"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.59%  "

$ws.Range("E48").Value = "  -2.11%  "

$ws.Range("D49").Value = "2.606.99"
$ws.Range("E49").Value = "  +11.09%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "6.96"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.49%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "23.39"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.20%  "
